# Applies two text edits (each splitting one run into several runs):
#   1) Slide 11, content placeholder, paragraph 3:
#        "细搜（主要成分，重心成分，带权评估）。"
#      -> "细搜（主要" / "成份，" / "重心" / "成份，" / "带权评估）。"
#         (note: 成分 -> 成份, twice)
#   2) Slide 7, content placeholder, paragraph 3:
#        "这个速度和精度，保证了德塔养料经华瑞集的极速分析能力。"
#      -> "这个速度和精度，保证了德塔" / "养" / "疗" / "经" / "华瑞集的极速分析能力。"
#         (note: 养料经 -> 养疗经)

$p = $ppt.ActivePresentation

# ---- Slide 11: "细搜（主要成分，重心成分，带权评估）。" ----
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$tr11 = $sh11.TextFrame.TextRange

$full11 = $tr11.Text
$idx11 = $full11.IndexOf("细搜（主要成分，重心成分，带权评估）。")
$start11 = $idx11 + 1

$pieces11 = @("细搜（主要", "成份，", "重心", "成份，", "带权评估）。")
$pos = $start11
foreach ($piece in $pieces11) {
    $len = $piece.Length
    $tr11.Characters($pos, $len).Text = $piece
    $pos = $pos + $len
}

# ---- Slide 7: "这个速度和精度，保证了德塔养料经华瑞集的极速分析能力。" ----
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange

$full7 = $tr7.Text
$idx7 = $full7.IndexOf("这个速度和精度，保证了德塔养料经华瑞集的极速分析能力。")
$start7 = $idx7 + 1

$pieces7 = @("这个速度和精度，保证了德塔", "养", "疗", "经", "华瑞集的极速分析能力。")
$pos = $start7
foreach ($piece in $pieces7) {
    $len = $piece.Length
    $tr7.Characters($pos, $len).Text = $piece
    $pos = $pos + $len
}
